$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$sourceValue = $ws.Range("A1679").Value2

for ($r = 1680; $r -le 1751; $r++) {
    $ws.Cells.Item($r, 1).Value2 = $sourceValue
}
